# Auto-Eval-ASI-D1.xlsx — fill in "Thomas Perreyon" grade column (N) with
# the evaluation results and update the saved view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Auto-Evaluation")

# Grades entered in column N (per-criterion evaluation for the 4th reviewer)
$ws.Range("N7").Value  = "C"
$ws.Range("N8").Value  = "C"
$ws.Range("N9").Value  = "B"
$ws.Range("N10").Value = "C"
$ws.Range("N11").Value = "B"
$ws.Range("N13").Value = "B"
$ws.Range("N14").Value = "B"
$ws.Range("N15").Value = "B"
$ws.Range("N16").Value = "B"
$ws.Range("N17").Value = "B"
$ws.Range("N18").Value = "B"
$ws.Range("N19").Value = "B"
$ws.Range("N20").Value = "C"
$ws.Range("N21").Value = "C"
$ws.Range("N22").Value = "B"
$ws.Range("N23").Value = "C"
$ws.Range("N24").Value = "B"
$ws.Range("N25").Value = "B"
$ws.Range("N26").Value = "B"
$ws.Range("N35").Value = "AB"
$ws.Range("N37").Value = "B"
$ws.Range("N42").Value = "A"
$ws.Range("N43").Value = "AB"
$ws.Range("N44").Value = "C"
$ws.Range("N46").Value = "AB"

# Restore the view: no frozen scroll position, selection parked on M48
$ws.Range("M48").Select()
